$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "1.000" or
# "4.383" are preserved verbatim instead of being coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "25.394.73"
$ws.Range("E2").Value = "  -1.08%  "
$ws.Range("D3").Value = "1.663.84"
$ws.Range("E3").Value = "  -1.65%  "
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  -0.60%  "
$ws.Range("D5").Value = "235.92"
$ws.Range("E5").Value = "  -2.20%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.43%  "
$ws.Range("D7").Value = "0.4776"
$ws.Range("E7").Value = "  -1.96%  "
$ws.Range("D8").Value = "0.2606"
$ws.Range("E8").Value = "  -1.85%  "
$ws.Range("D9").Value = "0.06147"
$ws.Range("E9").Value = "  +1.65%  "
$ws.Range("D10").Value = "0.07072"
$ws.Range("E10").Value = "  -1.37%  "
$ws.Range("D11").Value = "1.659.90"
$ws.Range("E11").Value = "  -3.28%  "
$ws.Range("D12").Value = "14.74"
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("D13").Value = "0.5918"
$ws.Range("E13").Value = "  -6.49%  "
$ws.Range("D14").Value = "4.383"
$ws.Range("E14").Value = "  -5.75%  "
$ws.Range("D15").Value = "74.42"
$ws.Range("E15").Value = "  -0.35%  "
$ws.Range("D16").Value = "0.9998"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("D17").Value = "1.000"
$ws.Range("E17").Value = "  -0.60%  "
$ws.Range("D18").Value = "25.400.21"
$ws.Range("E18").Value = "  -1.07%  "
$ws.Range("D19").Value = "0.000006766"
$ws.Range("E19").Value = "  +1.16%  "
$ws.Range("D20").Value = "11.41"
$ws.Range("E20").Value = "  -1.48%  "
$ws.Range("D21").Value = "1.873.08"
$ws.Range("E21").Value = "  -2.80%  "
$ws.Range("D22").Value = "4.439"
$ws.Range("E22").Value = "  -1.01%  "
$ws.Range("D23").Value = "8.651"
$ws.Range("E23").Value = "  -0.22%  "
$ws.Range("D24").Value = "5.340"
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").Value = "133.63"
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").Value = "15.05"
$ws.Range("E26").Value = "  +0.64%  "
$ws.Range("D27").Value = "1.397"
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").Value = "104.21"
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("D29").Value = "1.690"
$ws.Range("E29").Value = "  -2.66%  "
$ws.Range("D30").Value = "3.990"
$ws.Range("E30").Value = "  +3.35%  "
$ws.Range("D31").Value = "3.614"
$ws.Range("E31").Value = "  +1.02%  "
$ws.Range("D32").Value = "0.07641"
$ws.Range("E32").Value = "  -4.65%  "
$ws.Range("D33").Value = "0.04377"
$ws.Range("E33").Value = "  -5.54%  "
$ws.Range("D34").Value = "0.9995"
$ws.Range("E34").Value = "  -0.44%  "
$ws.Range("D35").Value = "2.604"
$ws.Range("E35").Value = "  -2.54%  "
$ws.Range("D36").Value = "0.6125"
$ws.Range("E36").Value = "  +3.60%  "
$ws.Range("D37").Value = "0.9440"
$ws.Range("E37").Value = "  -2.45%  "
$ws.Range("D38").Value = "2.621"
$ws.Range("E38").Value = "  -2.37%  "
$ws.Range("D39").Value = "0.8547"
$ws.Range("E39").Value = "  +0.95%  "
$ws.Range("D40").Value = "0.9997"
$ws.Range("E40").Value = "  -0.59%  "
$ws.Range("D41").Value = "0.01505"
$ws.Range("E41").Value = "  -4.20%  "
$ws.Range("D42").Value = "1.830"
$ws.Range("E42").Value = "  -3.30%  "
$ws.Range("D43").Value = "98.25"
$ws.Range("E43").Value = "  -1.43%  "
$ws.Range("D44").Value = "0.3766"
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("D45").Value = "4.666"
$ws.Range("E45").Value = "  -5.29%  "
$ws.Range("D46").Value = "0.1117"
$ws.Range("E46").Value = "  -3.42%  "
$ws.Range("D47").Value = "6.210"
$ws.Range("E47").Value = "  +1.02%  "
$ws.Range("D48").Value = "0.05256"
$ws.Range("E48").Value = "  +0.59%  "
$ws.Range("D49").Value = "29.55"
$ws.Range("E49").Value = "  -1.23%  "
$ws.Range("B50").Value = "TrueUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd"
$ws.Range("D50").Value = "1.001"
$ws.Range("E50").Value = "  -0.55%  "
$ws.Range("D51").Value = "0.3345"
$ws.Range("E51").Value = "  -1.13%  "
